$d = $word.ActiveDocument
$docXml = '<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"><w:body><w:p w14:paraId="0AF1CE24" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="343434"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="343434"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ideas for the Survival Game</w:t></w:r></w:p><w:p w14:paraId="797DE662" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="6DFC3D5A" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>AI:</w:t></w:r></w:p><w:p w14:paraId="00CA9A5D" w14:textId="0A020331" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="240"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Different countries and nations also have the behaviours tab and allow communications throw the different people this means that the way people communicate between each other could affect the relationship.</w:t></w:r></w:p><w:p w14:paraId="77908E37" w14:textId="7F4B7548" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Trading with other AI, other solar systems have a different set of rules so if entering the players entering must follow and respect the rules of the AI or the other players controlling the land otherwise punishment is under their command.</w:t></w:r></w:p><w:p w14:paraId="7FACC92A" w14:textId="6C5545DC" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Allow the AI nations to develop alongside of the players. </w:t></w:r></w:p><w:p w14:paraId="22AE5F4C" w14:textId="4A1D1218" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">AI have jobs within the nations that they are part of allowing the nations to function properly. </w:t></w:r></w:p><w:p w14:paraId="168AB703" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="629B9DBE" w14:textId="089DF879" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="007C1F36"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>pace</w:t></w:r><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p w14:paraId="06FF6FF2" w14:textId="78275478" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="240"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Dangers outside of the planets. e.g. this could be an incident within the rocket such as running out of fuel or colliding with an asteroid. Meaning the rocket would be stranded or damaging depending on how hard, fast and big the asteroid or waste material hit the rocket.</w:t></w:r></w:p><w:p w14:paraId="10765C4D" w14:textId="1931F6B1" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="240"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Asteroids, must implemented into the game</w:t></w:r></w:p><w:p w14:paraId="27BFB118" w14:textId="00E091A6" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="240"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Gas Giants, mining for power.</w:t></w:r></w:p><w:p w14:paraId="7FF9A541" w14:textId="3591EF60" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Ice giants to make oxygen and hydrogen. </w:t></w:r></w:p><w:p w14:paraId="2745B004" w14:textId="19FA53D0" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Black holes</w:t></w:r></w:p><w:p w14:paraId="377E579B" w14:textId="23628BD6" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Worm holes, ends up in a different part of the galaxy (corrupts and destroys ship)</w:t></w:r></w:p><w:p w14:paraId="04F39A19" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="106B5AD2" w14:textId="67E6A652" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>W</w:t></w:r><w:r w:rsidR="007C1F36"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>orld</w:t></w:r><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p w14:paraId="0187E470" w14:textId="5EEDDA7F" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Start from the beginning?</w:t></w:r></w:p><w:p w14:paraId="55534E5D" w14:textId="7A6EE906" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Have starter kit, e.g. farming tools and mining tools as well. General basic ideas and materials needed for the people to comfortably start the game.</w:t></w:r></w:p><w:p w14:paraId="77F23008" w14:textId="573299B2" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>What sort of obstacles can we throw in the way of the player to disrupt the difficulty level? Aggressiveness of animals and other creatures</w:t></w:r></w:p><w:p w14:paraId="21299DF1" w14:textId="505986A9" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Natural disasters on earth such as eruptions and earthquakes. This all depends on where in the world or plant you are on would cause the problem.  </w:t></w:r></w:p><w:p w14:paraId="549CE3EB" w14:textId="183AF758" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Different environments, extreme environments such as Deserts, or Antarctica, therefore only certain animals and organisms can live in these harsh </w:t></w:r><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>environments. Players would have to be careful with what they mess with o</w:t></w:r><w:r w:rsidR="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">r </w:t></w:r><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>how the approach the environment.</w:t></w:r></w:p><w:p w14:paraId="21E07F25" w14:textId="579935AC" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Adding seasons, Spring, Summer, Autumn &amp; Winter. Different planets may start at different season and some may last for less amounts of time or possibly longer. </w:t></w:r></w:p><w:p w14:paraId="3C7A80F0" w14:textId="5138A8BB" w:rsidR="00AE0486" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>People divided up between different countries, maybe implement the current world countries and conflicts. Change the situations like the world every update?</w:t></w:r></w:p><w:p w14:paraId="05ADEDB2" w14:textId="21D8F4AB" w:rsidR="00CF3070" w:rsidRDefault="00CF3070" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Where do we all start?</w:t></w:r></w:p><w:p w14:paraId="1B2F13EC" w14:textId="5BB5B5A0" w:rsidR="00CF3070" w:rsidRPr="005C4A65" w:rsidRDefault="00CF3070" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Start at a space fight, a scenario could be you’re in a dog fight and the ship is shot down, you manage to survive, and you are a planet that you don’t know of. (Earth), you find people and begin life the way they are living it. In the end you develop and become intelligent enough make your way back to where you came from.</w:t></w:r></w:p><w:p w14:paraId="42B15081" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="276A2C1F" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Resources:</w:t></w:r></w:p><w:p w14:paraId="2F2B7123" w14:textId="084B3557" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Ores e.g. Iron, Gold &amp; Copper. </w:t></w:r></w:p><w:p w14:paraId="30A01752" w14:textId="5A26977B" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Alloys, Bronze and other materials that are made from other metals. Use a smeltery or some form of smelting system to use the ores.</w:t></w:r></w:p><w:p w14:paraId="08AB62D5" w14:textId="09D8B323" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Basic minerals such as Iron has a symbol Fe, and the atom value of such and such… if we implement the periodic table and all the elements with our own as well the different compounds, mixtures and minerals we could make would be awesome.</w:t></w:r></w:p><w:p w14:paraId="31FB4824" w14:textId="61B35D3C" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Fuels, fossil fuels and other types of fuels would be a really good idea, because the players must power the bases, vehicles and other entities. </w:t></w:r></w:p><w:p w14:paraId="100B60DA" w14:textId="36ABE5AB" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Those fuels have to be refined and fractions in way so the players would have to find a way how to do this.</w:t></w:r></w:p><w:p w14:paraId="1F8F4185" w14:textId="5E65CE52" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Renewable power, such as wind power or hydroelectric. This could be done on a smaller scale with a little river or the players could use the initiative to make the power or design a more effective way.</w:t></w:r></w:p><w:p w14:paraId="55D1593E" w14:textId="2BB63689" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Renewable resources such as wood and other plants that grow. For example, at a really simple stage you could use the fibres from stinging nettles. </w:t></w:r></w:p><w:p w14:paraId="5B5975C4" w14:textId="0A48C6FE" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Resources such as trees would have to regrow, what about following the biodiversity and the regrowth of plants. E.g. bees pollenating flowers and birds taking and processing seeds</w:t></w:r></w:p><w:p w14:paraId="1207BCAB" w14:textId="437C16A8" w:rsidR="00AE0486" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Food chain? Microorganisms - little fish - larger fish - etc.</w:t></w:r></w:p><w:p w14:paraId="0A7E155C" w14:textId="3E2CC224" w:rsidR="000C218D" w:rsidRDefault="000C218D" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Grades of an ore could be a late game factor, for example the bioleaching could be used to improve </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> ore, if this was done on an industrious scale. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Profits would be incredibly good for the people who developed it and other benefits would come for them as well as good quality ores.</w:t></w:r></w:p><w:p w14:paraId="51CC07FE" w14:textId="77777777" w:rsidR="000C218D" w:rsidRPr="000C218D" w:rsidRDefault="000C218D" w:rsidP="000C218D"><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="20"/><w:tab w:val="left" w:pos="261"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="4DDDA524" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="4B1538AF" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Players: </w:t></w:r></w:p><w:p w14:paraId="082C6CEA" w14:textId="216A10F9" w:rsidR="00AE0486" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Players earn money by having a job and a salary. This job could be anything you could dream of as longs the job is offered by the nations. Remember different nations may have different jobs because of their situation and the circumstances they live in.</w:t></w:r></w:p><w:p w14:paraId="7D19E816" w14:textId="3101616D" w:rsidR="00F71D8C" w:rsidRPr="005C4A65" w:rsidRDefault="00F71D8C" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Jobs pay, this could be functional in the real-life way. Employed by large companies and work on behalf of them or start you own business competing against others and sell you skills and materials to other or people to make money. Depending on how well you work is dependent on how much you receive in return. </w:t></w:r></w:p><w:p w14:paraId="2B5B1098" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Hosting:</w:t></w:r></w:p><w:p w14:paraId="20B7D795" w14:textId="4636B8C1" w:rsidR="00737283" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">8 players at most, a recommended amount would be 4. The more players the worse the hosting becomes for the person hosting and the players. </w:t></w:r></w:p><w:p w14:paraId="754A68B9" w14:textId="5FE9A04F" w:rsidR="00AE0486" w:rsidRPr="005C4A65" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Hosting via servers could be a massive potential, look at GTA servers and other servers similar to them. This could be a potential way forward. Play hosting would be first before probably though. </w:t></w:r></w:p><w:p w14:paraId="42A0799E" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="00AE0486"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="21667858" w14:textId="77777777" w:rsidR="00AE0486" w:rsidRPr="00737283" w:rsidRDefault="008B2BCC"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00737283"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Laws of the Land:</w:t></w:r></w:p><w:p w14:paraId="7323888E" w14:textId="2F51F36E" w:rsidR="008B2BCC" w:rsidRDefault="008B2BCC" w:rsidP="005C4A65"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">All well developed or developing nations should have a set of rules that lay out what should </w:t></w:r><w:r w:rsidR="00737283" w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>happen,</w:t></w:r><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and players/AI should follow these rules otherwise punishments would be enforced. For </w:t></w:r><w:r w:rsidR="00737283" w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>example,</w:t></w:r><w:r w:rsidRPr="005C4A65"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the ‘major nations’ would council if the crime committed was incredibly bad. Death punishments may be enforced elsewhere as supposed to different ones.</w:t></w:r></w:p><w:p w14:paraId="75F482A6" w14:textId="556D9440" w:rsidR="007C1F36" w:rsidRDefault="007C1F36" w:rsidP="007C1F36"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Game Modes:</w:t></w:r></w:p><w:p w14:paraId="56BB5F1E" w14:textId="77777777" w:rsidR="007C1F36" w:rsidRDefault="007C1F36" w:rsidP="007C1F36"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Dog fighting game mode, this is the starting entry point to the game whilst other features are developed around the game. This game would be a multiplayer game where battles take place around the galaxy. The obstacles could be flying bits of debris from asteroids, destroyed ships &amp; other objects similar to the scene. </w:t></w:r></w:p><w:p w14:paraId="67DD4413" w14:textId="0B5F6E5B" w:rsidR="007C1F36" w:rsidRDefault="007C1F36" w:rsidP="007C1F36"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Tests ship flight controls, and other features such as the control panel and the field of view within the plane: 1</w:t></w:r><w:r w:rsidRPr="007C1F36"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> person, 3</w:t></w:r><w:r w:rsidRPr="007C1F36"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> person, look down view. These </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">are all ideas that could make the players experience better for the game and the game mode.  </w:t></w:r></w:p><w:p w14:paraId="6719468A" w14:textId="180AF3A8" w:rsidR="007139BF" w:rsidRDefault="00881523" w:rsidP="007C1F36"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Heists these could consist on taking out a particular target e.g. a nation leader that is corrupt and you have granted access with you and your team to take out this supposed target. Another one could be to bo</w:t></w:r><w:r w:rsidR="00F71D8C"><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>d a vessel and take the contraband kill enemies if necessary.</w:t></w:r></w:p><w:p w14:paraId="7FDF30E8" w14:textId="61DAC86A" w:rsidR="00F71D8C" w:rsidRDefault="00F71D8C" w:rsidP="007C1F36"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Capture the flag, for instance whilst the game is in development for larger parts and aspects of the game. We could implement a game mode alongside of the dog fighting too see what players like about the system and how the game could be improved or particular parts about the game.</w:t></w:r></w:p><w:p w14:paraId="21476925" w14:textId="4B381948" w:rsidR="00F71D8C" w:rsidRDefault="00F71D8C" w:rsidP="00F71D8C"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Economy:</w:t></w:r></w:p><w:p w14:paraId="422EA9A3" w14:textId="0916886D" w:rsidR="00F71D8C" w:rsidRPr="00F71D8C" w:rsidRDefault="00F71D8C" w:rsidP="00F71D8C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>All countries would need a form of trade, economy or way to thrive otherwise it</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>’</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s going to be difficult to survive in the harsh reality of the world. Other nations will be more power and successful. Players could follow in suit and sell more items and need resource to climb the graph.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Countries must make a stable way to control this system otherwise they could fail to succeed and come into massive problems. As the more money you pump into the system the less the money is worth. Therefore, creating a special way for the AI to control it would be easiest and the most efficient method.</w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Jobs:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="288" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Does your character have the key roles to work in this job? Is the intellectual of the character high enough. What sort of pay is there? Is this person reputable. All these things make you the person for the job. Making you intellect, strength and other aspects could be implemented to make the game more difficult and more realistic for the players. As this could be a proper experience.</w:t></w:r></w:p><w:sectPr w:rsidR="00F71D8C" w:rsidRPr="00F71D8C"><w:pgSz w:w="11900" w:h="16840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:noEndnote/></w:sectPr></w:body></w:document>'
$numberingXml = '<w:numbering xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"><w:abstractNum w:abstractNumId="0" w15:restartNumberingAfterBreak="0"><w:nsid w:val="00000001"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="00000001"/><w:lvl w:ilvl="0" w:tplc="00000001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="•"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="2" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="3" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="4" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="5" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="6" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="7" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="8" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="00000002"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="00000002"/><w:lvl w:ilvl="0" w:tplc="00000065"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="•"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="2" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="3" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="4" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="5" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="6" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="7" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="8" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="00000003"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="00000003"/><w:lvl w:ilvl="0" w:tplc="000000C9"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="•"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="2" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="3" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="4" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="5" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="6" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="7" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="8" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="00000004"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="00000004"/><w:lvl w:ilvl="0" w:tplc="0000012D"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="•"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="2" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="3" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="4" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="5" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="6" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="7" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="8" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0"><w:nsid w:val="00000005"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="00000005"/><w:lvl w:ilvl="0" w:tplc="00000191"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="•"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="1" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="2" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="3" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="4" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="5" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="6" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="7" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl><w:lvl w:ilvl="8" w:tplc="FFFFFFFF"><w:numFmt w:val="decimal"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="5" w15:restartNumberingAfterBreak="0"><w:nsid w:val="28917727"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="E0E0B588"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="6" w15:restartNumberingAfterBreak="0"><w:nsid w:val="2EC03643"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="A81CC27E"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="7" w15:restartNumberingAfterBreak="0"><w:nsid w:val="45F16D95"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="3BAA6206"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="8" w15:restartNumberingAfterBreak="0"><w:nsid w:val="4F0A78FF"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="9FD2AC36"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="9" w15:restartNumberingAfterBreak="0"><w:nsid w:val="51BF58F5"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="0F824EC4"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="10" w15:restartNumberingAfterBreak="0"><w:nsid w:val="5BC565FB"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="F1D05614"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="11" w15:restartNumberingAfterBreak="0"><w:nsid w:val="6BAD2D7B"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="E9D08FB0"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="12" w15:restartNumberingAfterBreak="0"><w:nsid w:val="705074A8"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="2E2CA8A0"/><w:lvl w:ilvl="0" w:tplc="102CDAF2"><w:numFmt w:val="bullet"/><w:lvlText w:val="-"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Helvetica" w:cs="Helvetica" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="13" w15:restartNumberingAfterBreak="0"><w:nsid w:val="774440A7"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="82321C84"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:cs="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="14" w15:restartNumberingAfterBreak="0"><w:nsid w:val="7D1B4222"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="9656D7BE"/><w:lvl w:ilvl="0" w:tplc="08090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:cs="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="08090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="08090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="08090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:num w:numId="1"><w:abstractNumId w:val="0"/></w:num><w:num w:numId="2"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="3"><w:abstractNumId w:val="2"/></w:num><w:num w:numId="4"><w:abstractNumId w:val="3"/></w:num><w:num w:numId="5"><w:abstractNumId w:val="4"/></w:num><w:num w:numId="6"><w:abstractNumId w:val="10"/></w:num><w:num w:numId="7"><w:abstractNumId w:val="7"/></w:num><w:num w:numId="8"><w:abstractNumId w:val="5"/></w:num><w:num w:numId="9"><w:abstractNumId w:val="8"/></w:num><w:num w:numId="10"><w:abstractNumId w:val="9"/></w:num><w:num w:numId="11"><w:abstractNumId w:val="6"/></w:num><w:num w:numId="12"><w:abstractNumId w:val="12"/></w:num><w:num w:numId="13"><w:abstractNumId w:val="13"/></w:num><w:num w:numId="14"><w:abstractNumId w:val="11"/></w:num><w:num w:numId="15"><w:abstractNumId w:val="14"/></w:num></w:numbering>'
$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
$pkg += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $docXml + '</pkg:xmlData></pkg:part>'
$pkg += '<pkg:part pkg:name="/word/numbering.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.numbering+xml"><pkg:xmlData>' + $numberingXml + '</pkg:xmlData></pkg:part>'
$pkg += '</pkg:package>'
$r = $d.Content
$r.InsertXML($pkg)
Write-Host ("Paragraphs: " + $d.Paragraphs.Count)
